$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '30.385.91'
$ws.Range('E2').Value = '  +1.12%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.014.14'
$ws.Range('E3').Value = '  +5.18%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '325.16'
$ws.Range('E5').Value = '  +1.48%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '1.000'
$ws.Range('E6').Value = '  +0.02%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.5104'
$ws.Range('E7').Value = '  +1.21%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.4157'
$ws.Range('E8').Value = '  +3.48%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.08777'
$ws.Range('E9').Value = '  +6.41%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '1.136'
$ws.Range('E10').Value = '  +2.62%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '42.88'
$ws.Range('E11').Value = '  +1.96%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '24.65'
$ws.Range('E12').Value = '  +2.42%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '2.010.09'
$ws.Range('E13').Value = '  +4.96%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '6.606'
$ws.Range('E14').Value = '  +2.76%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '7.489'
$ws.Range('E15').Value = '  +2.78%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '1.003'
$ws.Range('E16').Value = '  +0.15%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '94.40'
$ws.Range('E17').Value = '  +2.49%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.00001116'
$ws.Range('E18').Value = '  +1.77%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.06529'
$ws.Range('E19').Value = '  +0.41%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '19.00'
$ws.Range('E20').Value = '  +5.21%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '1.000'
$ws.Range('E21').Value = '  +0.03%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.240'
$ws.Range('E22').Value = '  +4.94%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '30.438.88'
$ws.Range('E23').Value = '  +1.17%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '11.96'
$ws.Range('E24').Value = '  +5.82%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.226'
$ws.Range('E25').Value = '  +1.24%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.237.47'
$ws.Range('E26').Value = '  +4.77%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '22.37'
$ws.Range('E27').Value = '  -0.43%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '162.99'
$ws.Range('E28').Value = '  +0.76%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.423'
$ws.Range('E29').Value = '  +6.11%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '131.57'
$ws.Range('E30').Value = '  +1.84%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.140'
$ws.Range('E31').Value = '  +1.23%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.1053'
$ws.Range('E32').Value = '  +1.40%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '6.117'
$ws.Range('E33').Value = '  +2.04%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '3.828'
$ws.Range('E34').Value = '  +1.43%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.352'
$ws.Range('E35').Value = '  +12.94%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.02522'
$ws.Range('E36').Value = '  +3.23%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '5.473'
$ws.Range('E37').Value = '  +2.35%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.06620'
$ws.Range('E38').Value = '  +3.19%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '12.38'
$ws.Range('E39').Value = '  +8.50%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '9.137'
$ws.Range('E40').Value = '  +4.42%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.2197'
$ws.Range('E41').Value = '  +1.63%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.6684'
$ws.Range('E42').Value = '  +2.11%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.232'
$ws.Range('E43').Value = '  +0.93%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '13.58'
$ws.Range('E44').Value = '  +1.41%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.6186'
$ws.Range('E45').Value = '  +2.53%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.201'
$ws.Range('E46').Value = '  +0.12%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '3.672'
$ws.Range('E47').Value = '  +0.92%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.273'
$ws.Range('E48').Value = '  +4.84%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '124.37'
$ws.Range('E49').Value = '  +0.76%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '81.37'
$ws.Range('E50').Value = '  +3.52%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.06895'
$ws.Range('E51').Value = '  +1.49%  '
